$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 0.1723813810668744
$ws.Range("C2").Value = 0.3499908675985048
$ws.Range("D2").Value = 0.6691348296765106
$ws.Range("E2").Value = 0.8300803312953138
$ws.Range("F2").Value = 0.5122925844322204
$ws.Range("G2").Value = 0.3934798787210216
$ws.Range("H2").Value = 0.6635913754756477
$ws.Range("I2").Value = 0.5073009638889382
$ws.Range("J2").Value = 0.6562794631732567
$ws.Range("K2").Value = 0.3761651213703216
$ws.Range("B3").Value = 0.3724188472315134
$ws.Range("C3").Value = 0.6717383692389594
$ws.Range("D3").Value = 0.8080126456412985
$ws.Range("E3").Value = 0.5063106041278653
$ws.Range("F3").Value = 0.3881366285386105
$ws.Range("G3").Value = 0.6532228421052815
$ws.Range("H3").Value = 0.4988236417653144
$ws.Range("I3").Value = 0.648507154866446
$ws.Range("J3").Value = 0.367541147837353
$ws.Range("K3").Value = 0.2069891273051463
$ws.Range("B4").Value = 0.6931629406961795
$ws.Range("C4").Value = 0.9281904944405108
$ws.Range("D4").Value = 0.4171162597140823
$ws.Range("E4").Value = 0.3584096246900327
$ws.Range("F4").Value = 0.6618489621351915
$ws.Range("G4").Value = 0.4728834683398792
$ws.Range("H4").Value = 0.6257285918382185
$ws.Range("I4").Value = 0.3539231126887352
$ws.Range("J4").Value = 0.1884683301199705
$ws.Range("K4").Value = 0.0637708202161355
$ws.Range("B5").Value = 0.8858698621201214
$ws.Range("C5").Value = 0.3879117858583491
$ws.Range("D5").Value = 0.3553890248851115
$ws.Range("E5").Value = 0.6429876785237111
$ws.Range("F5").Value = 0.4525523372746749
$ws.Range("G5").Value = 0.6111984364039706
$ws.Range("H5").Value = 0.3372820671964672
$ws.Range("I5").Value = 0.1708844846912203
$ws.Range("J5").Value = 0.04729693668609269
$ws.Range("K5").Value = 0.5643005784802194
$ws.Range("B6").Value = 0.7281095264762592
$ws.Range("C6").Value = 0.4307600866690486
$ws.Range("D6").Value = 0.451894392210712
$ws.Range("E6").Value = 0.4772462498614081
$ws.Range("F6").Value = 0.6205778399602666
$ws.Range("G6").Value = 0.2824712685850762
$ws.Range("H6").Value = 0.1519903938025505
$ws.Range("I6").Value = 0.03243293587104695
$ws.Range("J6").Value = 0.5355888678508487
$ws.Range("K6").Value = 0.3304428394308724
$ws.Range("B7").Value = 0.8814549586592991
$ws.Range("C7").Value = 0.4990050944230704
$ws.Range("D7").Value = 0.2369630190705999
$ws.Range("E7").Value = 0.6565229210261294
$ws.Range("F7").Value = 0.2804205111454833
$ws.Range("G7").Value = 0.0754221543665203
$ws.Range("H7").Value = 0.005033541495763572
$ws.Range("I7").Value = 0.5092998176531234
$ws.Range("J7").Value = 0.2870455487175398
$ws.Range("J8").ClearContents()
$ws.Range("B8").Value = 0.8113306813466088
$ws.Range("C8").Value = 0.3702517588061242
$ws.Range("D8").Value = 0.4774444992099378
$ws.Range("E8").Value = 0.3086147499160505
$ws.Range("F8").Value = 0.1114474655589138
$ws.Range("G8").Value = -0.03250802870845004
$ws.Range("H8").Value = 0.5037574785657621
$ws.Range("I8").Value = 0.2912208776562884
$ws.Range("I9").ClearContents()
$ws.Range("B9").Value = 0.6058422499139654
$ws.Range("C9").Value = 0.5621232400842474
$ws.Range("D9").Value = 0.1630143832823334
$ws.Range("E9").Value = 0.1214012084181564
$ws.Range("F9").Value = -0.01714285900890511
$ws.Range("G9").Value = 0.464909143608339
$ws.Range("H9").Value = 0.2753750686291025
$ws.Range("H10").ClearContents()
$ws.Range("B10").Value = 0.8731474265656303
$ws.Range("C10").Value = 0.2801031576935286
$ws.Range("D10").Value = -0.04106293303141043
$ws.Range("E10").Value = 0.01181609145629647
$ws.Range("F10").Value = 0.5006520160632426
$ws.Range("G10").Value = 0.2440474222454754
$ws.Range("G11").ClearContents()
$ws.Range("B11").Value = 0.526928335097643
$ws.Range("C11").Value = -0.02340706549352092
$ws.Range("D11").Value = -0.08280614765469801
$ws.Range("E11").Value = 0.5328694043117227
$ws.Range("F11").Value = 0.2564355480731927
$ws.Range("F12").ClearContents()
$ws.Range("B12").Value = 0.216269293922691
$ws.Range("C12").Value = 0.002218095814515486
$ws.Range("D12").Value = 0.4163848259537086
$ws.Range("E12").Value = 0.2715408197250452
$ws.Range("E13").ClearContents()
$ws.Range("B13").Value = 0.1670781728486028
$ws.Range("C13").Value = 0.4299220982016248
$ws.Range("D13").Value = 0.2085679007350822
$ws.Range("D14").ClearContents()
$ws.Range("B14").Value = 0.6837711400670328
$ws.Range("C14").Value = 0.3078859509171186
$ws.Range("C15").ClearContents()
$ws.Range("B15").Value = 0.3519456421565676
$ws.Range("B16").ClearContents()
